# Auto-generated edit script: update cryptos price table (D) and volume (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new value is numeric-looking (e.g. "0.9999", "20.34").
# Excel's normal type-inference would silently convert these to Number on assignment
# (dropping the original text formatting, e.g. "10.20" -> 10.2). The source data is
# textual (t="inlineStr" in the sheet XML), so force the cells to Text first, same as
# a user pre-formatting a column as Text before typing numeric-looking strings into it.
# (Done one cell at a time -- multi-area "A1,A2,..." Range strings only touch the
# first area in this host, so a loop is used instead of a single combined Range.)
foreach ($r in @(4, 5, 6, 7, 8, 9, 10, 11, 12, 14, 15, 16, 18, 19, 21, 22, 23, 24, 25, 26, 27, 28, 29, 30, 31, 32, 33, 34, 35, 36, 37, 38, 39, 40, 41, 42, 43, 44, 45, 47, 48, 49, 50, 51)) {
    $ws.Cells.Item($r, 4).NumberFormat = "@"
}

$ws.Range("D2").Value = "27.037.04"
$ws.Range("E2").Value = "  +5.56%  "
$ws.Range("D3").Value = "1.880.10"
$ws.Range("D4").Value = "0.9999"
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").Value = "281.39"
$ws.Range("E5").Value = "  +2.49%  "
$ws.Range("D6").Value = "0.9999"
$ws.Range("E6").Value = "  -0.16%  "
$ws.Range("D7").Value = "0.5264"
$ws.Range("E7").Value = "  +4.84%  "
$ws.Range("D8").Value = "0.3527"
$ws.Range("E8").Value = "  +0.95%  "
$ws.Range("D9").Value = "0.07042"
$ws.Range("E9").Value = "  +6.75%  "
$ws.Range("D10").Value = "20.34"
$ws.Range("E10").Value = "  +2.21%  "
$ws.Range("D11").Value = "0.8170"
$ws.Range("E11").Value = "  -1.95%  "
$ws.Range("D12").Value = "0.07795"
$ws.Range("E12").Value = "  +0.38%  "
$ws.Range("D13").Value = "1.879.77"
$ws.Range("E13").Value = "  +4.08%  "
$ws.Range("D14").Value = "5.218"
$ws.Range("E14").Value = "  +3.21%  "
$ws.Range("D15").Value = "90.61"
$ws.Range("E15").Value = "  +3.80%  "
$ws.Range("D16").Value = "0.9993"
$ws.Range("E16").Value = "  -0.13%  "
$ws.Range("E17").Value = "  +5.33%  "
$ws.Range("D18").Value = "0.000008159"
$ws.Range("E18").Value = "  +2.76%  "
$ws.Range("D19").Value = "1.000"
$ws.Range("E19").Value = "  -0.07%  "
$ws.Range("D20").Value = "27.064.11"
$ws.Range("E20").Value = "  +5.41%  "
$ws.Range("D21").Value = "4.768"
$ws.Range("E21").Value = "  +1.19%  "
$ws.Range("D22").Value = "10.20"
$ws.Range("E22").Value = "  +1.90%  "
$ws.Range("D23").Value = "6.236"
$ws.Range("E23").Value = "  +3.19%  "
$ws.Range("D24").Value = "2.396"
$ws.Range("E24").Value = "  +13.90%  "
$ws.Range("D25").Value = "146.65"
$ws.Range("E25").Value = "  +3.39%  "
$ws.Range("D26").Value = "17.60"
$ws.Range("E26").Value = "  +4.25%  "
$ws.Range("D27").Value = "1.675"
$ws.Range("E27").Value = "  +1.37%  "
$ws.Range("D28").Value = "113.54"
$ws.Range("E28").Value = "  +5.07%  "
$ws.Range("D29").Value = "4.389"
$ws.Range("E29").Value = "  +1.81%  "
$ws.Range("D30").Value = "4.384"
$ws.Range("E30").Value = "  +4.81%  "
$ws.Range("D31").Value = "0.08902"
$ws.Range("E31").Value = "  +1.46%  "
$ws.Range("D32").Value = "0.04915"
$ws.Range("E32").Value = "  +2.37%  "
$ws.Range("D33").Value = "1.177"
$ws.Range("E33").Value = "  +4.65%  "
$ws.Range("D34").Value = "0.7440"
$ws.Range("E34").Value = "  +3.40%  "
$ws.Range("D35").Value = "2.894"
$ws.Range("E35").Value = "  +0.48%  "
$ws.Range("D36").Value = "3.297"
$ws.Range("E36").Value = "  +9.00%  "
$ws.Range("D37").Value = "2.412"
$ws.Range("E37").Value = "  +6.22%  "
$ws.Range("D38").Value = "0.5307"
$ws.Range("E38").Value = "  +3.11%  "
$ws.Range("D39").Value = "0.01885"
$ws.Range("E39").Value = "  +1.46%  "
$ws.Range("D40").Value = "0.9836"
$ws.Range("E40").Value = "  +3.43%  "
$ws.Range("D41").Value = "117.11"
$ws.Range("E41").Value = "  +2.91%  "
$ws.Range("D42").Value = "6.322"
$ws.Range("E42").Value = "  +2.64%  "
$ws.Range("D43").Value = "8.185"
$ws.Range("E43").Value = "  +2.53%  "
$ws.Range("D44").Value = "0.9993"
$ws.Range("E44").Value = "  -0.14%  "
$ws.Range("D45").Value = "0.4607"
$ws.Range("E45").Value = "  +1.61%  "
$ws.Range("E46").Value = "  -0.62%  "
$ws.Range("D47").Value = "9.508"
$ws.Range("E47").Value = "  +2.35%  "
$ws.Range("D48").Value = "36.78"
$ws.Range("E48").Value = "  +2.66%  "
$ws.Range("D49").Value = "1.522"
$ws.Range("E49").Value = "  +2.38%  "
$ws.Range("D50").Value = "0.05947"
$ws.Range("E50").Value = "  +2.63%  "
$ws.Range("D51").Value = "61.87"
$ws.Range("E51").Value = "  +4.20%  "
